$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "edit1"
$ws.Range("B7").Value = "riya-morankar"
$ws.Range("C7").Value = "Merged"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2025-06-18"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "N/A"
